$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.409.54"
$ws.Range("E2").Value = "  +2.03%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.908.22"
$ws.Range("E3").Value = "  +0.78%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "528.89"
$ws.Range("E5").Value = "  +9.60%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.76"
$ws.Range("E6").Value = "  +0.47%  "
$ws.Range("E7").Value = "  -1.14%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.998"
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.720"
$ws.Range("E9").Value = "  -2.20%  "
$ws.Range("E10").Value = "  -3.56%  "
$ws.Range("E11").Value = "  -4.44%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "42.17"
$ws.Range("E12").Value = "  -1.27%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.534.60"
$ws.Range("E13").Value = "  +0.90%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.26"
$ws.Range("E14").Value = "  -2.54%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.954.88"
$ws.Range("E15").Value = "  +2.16%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.05"
$ws.Range("E16").Value = "  -1.53%  "
$ws.Range("E17").Value = "  +8.43%  "
$ws.Range("E18").Value = "  -0.50%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "19.75"
$ws.Range("E19").Value = "  -1.62%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "69.376.61"
$ws.Range("E20").Value = "  +1.94%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "430.17"
$ws.Range("E21").Value = "  +0.38%  "
$ws.Range("E22").Value = "  -4.70%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.21"
$ws.Range("E23").Value = "  -3.57%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "88.68"
$ws.Range("E24").Value = "  -1.49%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.04"
$ws.Range("E25").Value = "  +9.64%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.61"
$ws.Range("E26").Value = "  -2.88%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.61"
$ws.Range("E27").Value = "  -3.46%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "36.38"
$ws.Range("E28").Value = "  -2.62%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "685.26"
$ws.Range("E29").Value = "  -3.66%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "13.18"
$ws.Range("E30").Value = "  -2.24%  "
$ws.Range("E31").Value = "  -2.24%  "
$ws.Range("E32").Value = "  -2.89%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "67.83"
$ws.Range("E33").Value = "  +11.63%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.445"
$ws.Range("E34").Value = "  +12.62%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.01"
$ws.Range("E35").Value = "  -0.89%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "40.10"
$ws.Range("E36").Value = "  -1.90%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0₃0847"
$ws.Range("E37").Value = "  -2.75%  "
$ws.Range("E38").Value = "  +3.10%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.999"
$ws.Range("E39").Value = "  +0.07%  "
$ws.Range("E40").Value = "  -0.09%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0481"
$ws.Range("E41").Value = "  -3.67%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.16"
$ws.Range("E42").Value = "  +2.82%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.13"
$ws.Range("E43").Value = "  +5.48%  "
$ws.Range("E44").Value = "  -5.29%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.37"
$ws.Range("E45").Value = "  -0.43%  "
$ws.Range("E46").Value = "  -1.10%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0₆0359"
$ws.Range("E47").Value = "  +11.73%  "
$ws.Range("E48").Value = "  +7.61%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.745.83"
$ws.Range("E49").Value = "  +13.34%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "144.73"
$ws.Range("E50").Value = "  +0.03%  "
$ws.Range("E51").Value = "  -2.91%  "
